$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before row 2, pushing the existing data (old rows 2-21)
# down to rows 11-30. New rows 2-10 take on row 1's formatting by default,
# so clear that back off to match the plain (unstyled) data rows.
$ws.Range("A2:A10").EntireRow.Insert()
$ws.Range("A2:C10").ClearFormats()

# Full data block for rows 2-31 after the insertion (first 9 rows are the
# newly-added points, rows 11-30 are the original rows 2-21 shifted down
# by 9, and row 31 is a brand-new row appended at the end).
$data = @(
    @(-2.203232884407044, 3.565518021583557, 2.421210885047912),
    @(-2.176557183265686, 3.694440901279449, 2.333479553461074),
    @(-2.326471328735352, 3.527873754501343, 2.532370328903198),
    @(-2.219938099384308, 3.655792444944382, 2.356739014387131),
    @(-2.15459930896759, 3.624111294746399, 2.422300338745117),
    @(-2.086719453334808, 3.556297302246094, 2.323985010385513),
    @(-2.057324409484864, 3.517749786376953, 2.224772691726685),
    @(-2.064002573490143, 3.512480080127716, 2.171477824449539),
    @(-2.053778767585754, 3.516274869441986, 2.167593479156494),
    @(-2.025566756725311, 3.52062651515007, 2.27691987156868),
    @(-2.230706214904786, 3.561713695526123, 2.031704187393189),
    @(-2.475497364997865, 3.523229420185089, 2.537566900253297),
    @(-2.48767375946045, 3.061535835266112, 3.86157149076462),
    @(-2.442349374294281, 2.902286112308502, 4.137762367725372),
    @(-2.809046030044556, 1.600962877273557, 4.795935153961183),
    @(-3.426159977912904, -0.2882512211799662, 5.237384021282194),
    @(-4.615099787712106, 0.4410536289215286, 4.911977410316469),
    @(-5.570275843143422, 5.792195498943334, 7.507518291473446),
    @(4.738878250122034, 7.207733154296871, 21.9623451232909),
    @(-3.471131086349493, 6.420480489730835, -4.152171969413772),
    @(-2.357546925544736, 3.453876137733455, 3.317040443420424),
    @(-0.1309916377067535, 3.795689940452578, 3.876870155334465),
    @(-0.2815589904785198, 4.85230040550232, 2.205311059951784),
    @(-1.535586237907411, 4.95532476902008, 1.937351673841474),
    @(-1.377771139144897, 4.249351501464844, 1.058028712868691),
    @(0.1868795156478872, 3.352623224258432, 1.056536458432678),
    @(-0.04321670532226736, 5.682518005371086, 1.995282649993894),
    @(-0.4393689632415791, 3.92675977945327, 1.474137753248212),
    @(-0.5461759567260746, 4.010827064514163, 1.266485691070556),
    @(-0.5815373659133911, 4.265253961086274, 1.277935206890107)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
